# The workbook originally has three sheets:
#   "управление" (first, was not the active tab)
#   "факультет"  (second, was the active tab with selection C2:C9)
#   "Sheet1"     (third, empty, untouched)
#
# This edit:
#   1. Appends the English field-name hint to each of the 5 header cells
#      (row 1, columns A:E) on both the "управление" and "факультет"
#      sheets.
#   2. Swaps which sheet is active/selected: "управление" becomes the
#      active (tabSelected) sheet with A1:E1 selected, and "факультет"
#      loses tabSelected, also ending up with A1:E1 selected.

$wb = $excel.ActiveWorkbook

# Ordered (not a hashtable) so the shared-string table is rebuilt in a
# deterministic, reproducible sequence: A1, B1, C1, E1, D1.
$headerUpdates = @(
    @("A1", "тип категории (name)"),
    @("B1", "идентификаторы (id)"),
    @("C1", "актуальность (relevance)"),
    @("E1", "перевод (items)"),
    @("D1", "полное название элемента (keys)")
)

$sheetNames = @("управление", "факультет")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($pair in $headerUpdates) {
        $ws.Range($pair[0]).Value = $pair[1]
    }
}

# "факультет" is no longer the active tab, but still gets A1:E1 selected.
# Select it first (selecting implicitly activates a sheet), then activate
# and select "управление" last so it ends up as the final active tab.
$wsFak = $wb.Worksheets.Item("факультет")
$wsFak.Activate()
$wsFak.Range("A1:E1").Select()

# Make "управление" the active sheet/tab, with A1:E1 selected.
$wsUpr = $wb.Worksheets.Item("управление")
$wsUpr.Activate()
$wsUpr.Range("A1:E1").Select()
